$wb = $excel.ActiveWorkbook

# --- Update the Date value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-02-20T18:59:33+00:00"

# --- Add two new concept rows (JSON, HTML) to the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

# Row 7: JSON
$concepts.Range("A6:D6").Copy()
$concepts.Range("A7:D7").PasteSpecial(-4122)
$concepts.Range("A6").Copy()
$concepts.Range("A7").PasteSpecial(-4163)
$concepts.Range("B7").Value = "JSON"
$concepts.Range("C7").Value = "JSON File"

# Row 8: HTML
$concepts.Range("A6:D6").Copy()
$concepts.Range("A8:D8").PasteSpecial(-4122)
$concepts.Range("A6").Copy()
$concepts.Range("A8").PasteSpecial(-4163)
$concepts.Range("B8").Value = "HTML"
$concepts.Range("C8").Value = "HTML File"

$excel.CutCopyMode = 0
